# Fruta / hortaliza, semanal
# Insert two new weekly records at the top of the Chirimoya data block
# (rows 50-51), pushing the existing rows 50-77 down to 52-79.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 50:77 down by two rows -> 52:79
$ws.Rows("50:51").Insert()

# New row 50
$ws.Range("A50").Value = 10
$ws.Range("B50").Value = "Vega Modelo de Temuco"
$ws.Range("C50").Value = "La Araucanía"
$ws.Range("D50").Value = 44489
$ws.Range("E50").Value = 9
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100107
$ws.Range("H50").Value = "Otros"
$ws.Range("I50").Value = 100107002
$ws.Range("J50").Value = "Chirimoya"
$ws.Range("K50").Value = "Cultivar IV Región"
$ws.Range("L50").Value = "Especial"
$ws.Range("M50").Value = 30
$ws.Range("N50").Value = 2800
$ws.Range("O50").Value = 2800
$ws.Range("P50").Value = 2800
$ws.Range("Q50").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R50").Value = "Provincia del Elquí"
$ws.Range("S50").Value = 2800
$ws.Range("T50").Value = 1

# New row 51
$ws.Range("A51").Value = 10
$ws.Range("B51").Value = "Vega Modelo de Temuco"
$ws.Range("C51").Value = "La Araucanía"
$ws.Range("D51").Value = 44489
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100107
$ws.Range("H51").Value = "Otros"
$ws.Range("I51").Value = 100107002
$ws.Range("J51").Value = "Chirimoya"
$ws.Range("K51").Value = "Cultivar IV Región"
$ws.Range("L51").Value = "Primera"
$ws.Range("M51").Value = 110
$ws.Range("N51").Value = 2500
$ws.Range("O51").Value = 3500
$ws.Range("P51").Value = 2773
$ws.Range("Q51").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R51").Value = "Provincia del Elquí"
$ws.Range("S51").Value = 2773
$ws.Range("T51").Value = 1
